# daily auto push: 2026-01-24 06:47 UTC
# A new observation for 2026/01/24 (Sat) at hour 13 was logged before the
# 2026/12/29 entry, so every row from the old row 703 onward shifts down
# by one. Insert a fresh row at 703 and populate it; Excel's Insert()
# takes care of pushing rows 703..744 down to 704..745 (and the sheet's
# used-range dimension grows from D744 to D745 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(703).Insert()

# Column A holds plain text dates ("2026/01/24"), not real Excel dates.
# Force text storage first so the COM Value setter doesn't auto-coerce
# the slash-separated string into a date serial number, then drop the
# cell back to the same (unstyled) look as its neighbours.
$ws.Range("A703").NumberFormat = "@"
$ws.Range("A703").Value = "2026/01/24"
$ws.Range("A703").Style = $ws.Range("A702").Style

$ws.Range("B703").Value = "土"
$ws.Range("C703").Value = 13
$ws.Range("D703").Value = 201
